$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the contents of a given paragraph's two runs (the italic
# "speaker label" run and the normal "message" run) using a Find/Replace
# scoped strictly to that paragraph's Range, so we never touch unrelated
# occurrences of the same text elsewhere in the document.
# (NOTE: this runtime's PowerShell does not support named -Param syntax,
# so all functions here use positional parameters only.)
# ---------------------------------------------------------------------------
function Replace-InParagraph($ParaIndex, $OldLabel, $NewLabel, $OldMessage, $NewMessage) {
    $p = $d.Paragraphs.Item($ParaIndex)
    $rng1 = $p.Range
    $res1 = $rng1.Find.Execute($OldLabel, $false, $false, $false, $false, $false, $true, 0, $false, $NewLabel, 2)
    if (-not $res1) {
        throw "Label not found in paragraph $ParaIndex : $OldLabel"
    }

    $p2 = $d.Paragraphs.Item($ParaIndex)
    $rng2 = $p2.Range
    $res2 = $rng2.Find.Execute($OldMessage, $false, $false, $false, $false, $false, $true, 0, $false, $NewMessage, 2)
    if (-not $res2) {
        throw "Message not found in paragraph $ParaIndex : $OldMessage"
    }
}

# 1) "00:00 Bruno Torossi:" -> "00:00 Mariano Sonzini Astudillo:"
Replace-InParagraph 12 "00:00 Bruno Torossi:" "00:00 Mariano Sonzini Astudillo:" " Qué onda todo bien? qué onda boludo, cómo Va la vida." " Geográficas para que para que más gente se pueda se pueda sumar a ser parte o hacer usuario de lo que ofrece la tecnología cripto y blockchain."

# 2) "00:12 Gabriel E. Calvo:" -> "00:19 Bruno Torossi:"
Replace-InParagraph 13 "00:12 Gabriel E. Calvo:" "00:19 Bruno Torossi:" " vida de campo" " Cómo andan qué tal Nano hola Sandy"

# 3) "00:14 Bruno Torossi:" -> "00:20 Mariano Sonzini Astudillo:"
Replace-InParagraph 14 "00:14 Bruno Torossi:" "00:20 Mariano Sonzini Astudillo:" " Qué lindo boludo, qué onda está haciendo mucho frío?" " Muy bien."

# 4) "00:18 Gabriel E. Calvo:" -> "00:22 Bruno Torossi:"
Replace-InParagraph 15 "00:18 Gabriel E. Calvo:" "00:22 Bruno Torossi:" " mucho frío siempre sí, mucho frío un grado ponerle Y nada, una cabaña, media precaria digamos, o sea tiene, no sé precaria, soy abundante, pero tiene mucho vidrio mucho ventanal." " cómo va"

# 5) "00:32 Bruno Torossi:" -> "00:26 Mariano Sonzini Astudillo:"
Replace-InParagraph 16 "00:32 Bruno Torossi:" "00:26 Mariano Sonzini Astudillo:" " sí" " Estábamos con una breve intro Santi te parece que si vamos esperamos unos minutos, a ver si si aparece alguien."

# 6) "00:33 Gabriel E. Calvo:" -> "00:33 Santiago Cristobal:"
Replace-InParagraph 17 "00:33 Gabriel E. Calvo:" "00:33 Santiago Cristobal:" " Y el vidrio pasa todo el frío y solamente calefacción que tenemos no hay gas, no hay gas, no hay gas" " No me parece, me parece bien arrancando capaz de descubriendo gente capaz, que no en la idea es abierto y el que quiera se venga yo"

# 7) "00:41 Bruno Torossi (chat):" -> "00:39 Mariano Sonzini Astudillo:" / " más vale"
Replace-InParagraph 18 "00:41 Bruno Torossi (chat):" "00:39 Mariano Sonzini Astudillo:" " Hola, estoy transcribiendo esta llamada con mi extensión Tactiq AI: https://tactiq.io/r/transcribing" " más vale"

# ---------------------------------------------------------------------------
# Insert two brand-new transcript paragraphs right after paragraph 18
# (the one that now reads "00:39 Mariano Sonzini Astudillo: más vale").
# Each new paragraph mimics the same structure used throughout the
# transcript: an italic "speaker label" run followed by a plain message run.
# ---------------------------------------------------------------------------
function Insert-TranscriptParagraphAfter($AfterParaIndex, $Label, $Message) {
    $src = $d.Paragraphs.Item($AfterParaIndex)
    $src.Range.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($AfterParaIndex + 1)
    $labelStart = $newPara.Range.Start

    $insertPoint = $d.Range($labelStart, $labelStart)
    $insertPoint.InsertAfter($Label)
    $labelEnd = $labelStart + $Label.Length

    $msgPoint = $d.Range($labelEnd, $labelEnd)
    $msgPoint.InsertAfter($Message)

    $labelFmtRange = $d.Range($labelStart, $labelEnd)
    $labelFmtRange.Font.Italic = $true
}

Insert-TranscriptParagraphAfter 18 "00:41 Santiago Cristobal:" " lo quiero al baño, así que aprovecho"

Insert-TranscriptParagraphAfter 19 "00:42 Bruno Torossi (chat):" " Hola, estoy transcribiendo esta llamada con mi extensión Tactiq AI: https://tactiq.io/r/transcribing"

# ---------------------------------------------------------------------------
# Remove the two paragraphs that follow (originally "00:41 Gabriel E. Calvo:
# solamente es salamandra." and "00:43 Bruno Torossi: Claro, sí sí sí."),
# which are now paragraphs 21 and 22 after the two insertions above.
# ---------------------------------------------------------------------------
$delPara1 = $d.Paragraphs.Item(21)
$delPara1.Range.Text = ""
$delPara1.Range.Delete()

$delPara2 = $d.Paragraphs.Item(21)
$delPara2.Range.Text = ""
$delPara2.Range.Delete()

# ---------------------------------------------------------------------------
# Finally, the last transcript line ("00:46 Gabriel E. Calvo: Pero bueno,
# lindo hermoso...") becomes "00:42 Santiago Cristobal: que vos estás...".
# After the two deletions above it is paragraph 21 again.
# ---------------------------------------------------------------------------
Replace-InParagraph 21 "00:46 Gabriel E. Calvo:" "00:42 Santiago Cristobal:" " Pero bueno, lindo hermoso igual sea hermoso salió a correr lago otra vida." " que vos estás haciendo las introducciones para dejarlos un ratito por ahí se pueden presentar rápido o no hay grúas y vos también nos conoces a ellos y nada después si se va sumando gente genial."

Write-Host "Done. Final paragraph count:" $d.Paragraphs.Count
